$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")
$ws.Range("B3").Value = "Good"
$ws.Range("C3").Value = "Proactive part of the discord server, Motivated to work on the semi-solution for the server problems and machine learning research"
$ws.Range("B16").Value = "Good"
$ws.Range("C16").Value = "Good communication over Discord, Good commitments to GitHub and actively sharing files and experiences with other group members"

$ws.Range("C17").Select() | Out-Null
